$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank timesheet row 16 (10/13/2024, 2:04 PM - 2:54 PM)
$ws.Range("A16").Value = 45578
$ws.Range("B16").Value = 0.58611111111111114
$ws.Range("C16").Value = 0.62083333333333335

# Extend the weekly summary ranges in row 4 to include the new row 16 entry
$ws.Range("M4").Formula = "=SUM(D10:D16)"
$ws.Range("N4").Formula = "=SUM(G10:G16)"

# Update the active selection to reflect where the user ended up
$ws.Range("J14").Select()
